$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly cryptos data refresh (GitHub Actions style update).
# D-column price values that look numeric are forced to Text format
# first so Excel keeps them as literal strings (matching the source feed),
# e.g. "1.001" / "289.98" rather than being auto-converted to numbers.

$ws.Range("D2").Value = "22.016.07"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").Value = "1.552.63"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.001"
$ws.Range("E5").Value = "  -0.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "289.98"
$ws.Range("E6").Value = "  +0.75%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3957"
$ws.Range("E7").Value = "  +3.98%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3211"
$ws.Range("E8").Value = "  -2.98%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "43.85"
$ws.Range("E9").Value = "  -0.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07221"
$ws.Range("E10").Value = "  -2.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.073"
$ws.Range("E11").Value = "  -5.36%  "
$ws.Range("E12").Value = "  -0.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.679"
$ws.Range("E13").Value = "  -2.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.62"
$ws.Range("E14").Value = "  -6.89%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001126"
$ws.Range("E15").Value = "  +5.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.610"
$ws.Range("E16").Value = "  -1.51%  "
$ws.Range("D17").Value = "1.555.30"
$ws.Range("E17").Value = "  -5.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06588"
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "83.44"
$ws.Range("E19").Value = "  -2.90%  "
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.246"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.46"
$ws.Range("E22").Value = "  -3.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.28"
$ws.Range("E23").Value = "  -3.27%  "
$ws.Range("D24").Value = "22.027.77"
$ws.Range("E24").Value = "  -0.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.376"
$ws.Range("E25").Value = "  +3.78%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.405"
$ws.Range("E26").Value = "  -4.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "148.44"
$ws.Range("E27").Value = "  -1.42%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.53"
$ws.Range("E28").Value = "  -3.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.876"
$ws.Range("E29").Value = "  -1.06%  "
$ws.Range("D30").Value = "1.728.71"
$ws.Range("E30").Value = "  -4.81%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "118.62"
$ws.Range("E31").Value = "  -3.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9624"
$ws.Range("E32").Value = "  -11.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.819"
$ws.Range("E33").Value = "  -0.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08300"
$ws.Range("E34").Value = "  +1.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.154"
$ws.Range("E35").Value = "  -1.59%  "
$ws.Range("E36").Value = "  -16.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02255"
$ws.Range("E37").Value = "  -2.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.084"
$ws.Range("E38").Value = "  -3.53%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05981"
$ws.Range("E39").Value = "  -4.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.211"
$ws.Range("E40").Value = "  -1.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2029"
$ws.Range("E41").Value = "  -5.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.001"
$ws.Range("E42").Value = "  -0.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.69"
$ws.Range("E43").Value = "  -1.95%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5794"
$ws.Range("E44").Value = "  -4.07%  "
$ws.Range("B45").Value = "PancakeSwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.742"
$ws.Range("E45").Value = "  +0.17%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.96"
$ws.Range("E46").Value = "  -5.64%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5555"
$ws.Range("E47").Value = "  -4.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "118.03"
$ws.Range("E48").Value = "  -3.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.891"
$ws.Range("E49").Value = "  -3.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.131"
$ws.Range("E50").Value = "  -3.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06807"
$ws.Range("E51").Value = "  -3.10%  "
